$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 44539
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 3800
$ws.Range("O4").Value = 4000
$ws.Range("P4").Value = 3900
$ws.Range("R4").Value = 'Región del Maule'
$ws.Range("S4").Value = 1950

# Row 5
$ws.Range("D5").Value = 44174
$ws.Range("M5").Value = 150
$ws.Range("N5").Value = 3700
$ws.Range("O5").Value = 3800
$ws.Range("P5").Value = 3747
$ws.Range("S5").Value = 1874

# Row 6
$ws.Range("D6").Value = 44540
$ws.Range("M6").Value = 240
$ws.Range("N6").Value = 3500
$ws.Range("O6").Value = 3800
$ws.Range("P6").Value = 3650
$ws.Range("Q6").Value = '$/bandeja 2 kilos'
$ws.Range("R6").Value = 'Región del Maule'
$ws.Range("S6").Value = 1825
$ws.Range("T6").Value = 2

# Row 7
$ws.Range("D7").Value = 44187
$ws.Range("M7").Value = 80
$ws.Range("N7").Value = 2800
$ws.Range("O7").Value = 3000
$ws.Range("P7").Value = 2900
$ws.Range("R7").Value = 'Provincia de Linares'
$ws.Range("S7").Value = 1450

# Row 8
$ws.Range("D8").Value = 44187
$ws.Range("M8").Value = 65
$ws.Range("N8").Value = 1400
$ws.Range("O8").Value = 1500
$ws.Range("P8").Value = 1446
$ws.Range("Q8").Value = '$/envase 1 kilo'
$ws.Range("R8").Value = 'Provincia de Diguillín'
$ws.Range("S8").Value = 1446
$ws.Range("T8").Value = 1
